# Applies the "3e version avec organisation fichiers" edit to
# StructureDefinition-CapaciteSavoirfaire.xlsx:
#  - Metadata sheet: bump Date, change Base Definition URL
#  - Elements sheet: turn the old "capaciteSavoirFaire" row into the new
#    "typeSavoirFaire" row, insert two new date elements
#    (dateReconnaissance / dateAbandon), and append the (unchanged)
#    "capaciteSavoirFaire" row after them.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Metadata sheet
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-07-21T11:52:46+00:00"
$meta.Range("B18").Value = "https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/SavoirFaire"

# ---------------------------------------------------------------------
# 2) Elements sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Use row 3 (currently "CapaciteSavoirfaire.capaciteSavoirFaire") as the
# formatting template for the three new/edited data rows: copy its
# formats down into rows 4-6 before touching any values so every new
# cell keeps the same style as the existing data rows.
$ws.Range("A3:AJ3").Copy()
$ws.Range("A4:AJ6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 3: repurposed into "CapaciteSavoirfaire.typeSavoirFaire" -----
$ws.Range("A3").Value = "CapaciteSavoirfaire.typeSavoirFaire"
$ws.Range("B3").Value = "CapaciteSavoirfaire.typeSavoirFaire"
$ws.Range("L3").Value = " Le type de savoir-faire (qualifications/autres attributions) d" + [char]0x00E9 + "signe par exemple:** une sp" + [char]0x00E9 + "cialit" + [char]0x00E9 + " ordinale (S);** une comp" + [char]0x00E9 + "tence (C);** etc."
$ws.Range("M3").Value = $ws.Range("L3").Value
$ws.Range("Z3").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R04-TypeSavoirFaire?vs"
$ws.Range("AF3").Value = "SavoirFaire.typeSavoirFaire"
# F3/G3 (Min/Max), K3 (Type(s)), X3 (Binding Strength), AG3/AH3 (Base Min/Max)
# are unchanged ("0", "1", "Coding`n", "preferred", "0", "1").

# --- Row 4: new "CapaciteSavoirfaire.dateReconnaissance" --------------
$ws.Range("A4").Value = "CapaciteSavoirfaire.dateReconnaissance"
$ws.Range("B4").Value = "CapaciteSavoirfaire.dateReconnaissance"
$ws.Range("F4").Value = "0"
$ws.Range("G4").Value = "1"
$ws.Range("K4").Value = "date`n"
$ws.Range("L4").Value = " Date " + [char]0x00E0 + " laquelle, l" + [char]0x2019 + "organisme donnant l" + [char]0x2019 + "autorisation d" + [char]0x2019 + "exercer une qualification a reconnu cette qualification ou date " + [char]0x00E0 + " laquelle l'attribution a " + [char]0x00E9 + "t" + [char]0x00E9 + " donn" + [char]0x00E9 + "e au professionnel."
$ws.Range("M4").Value = $ws.Range("L4").Value
$ws.Range("X4").Value = ""
$ws.Range("Y4").Value = ""
$ws.Range("Z4").Value = ""
$ws.Range("AF4").Value = "SavoirFaire.dateReconnaissance"
$ws.Range("AG4").Value = "0"
$ws.Range("AH4").Value = "1"

# --- Row 5: new "CapaciteSavoirfaire.dateAbandon" ----------------------
$ws.Range("A5").Value = "CapaciteSavoirfaire.dateAbandon"
$ws.Range("B5").Value = "CapaciteSavoirfaire.dateAbandon"
$ws.Range("F5").Value = "0"
$ws.Range("G5").Value = "1"
$ws.Range("K5").Value = "date`n"
$ws.Range("L5").Value = " Date " + [char]0x00E0 + " laquelle le professionnel a d" + [char]0x00E9 + "clar" + [char]0x00E9 + " renoncer " + [char]0x00E0 + " l" + [char]0x2019 + "exercice d" + [char]0x2019 + "un savoir-faire ou date " + [char]0x00E0 + " laquelle il ne souhaite plus le faire appara" + [char]0x00EE + "tre."
$ws.Range("M5").Value = $ws.Range("L5").Value
$ws.Range("X5").Value = ""
$ws.Range("Y5").Value = ""
$ws.Range("Z5").Value = ""
$ws.Range("AF5").Value = "SavoirFaire.dateAbandon"
$ws.Range("AG5").Value = "0"
$ws.Range("AH5").Value = "1"

# --- Row 6: the original "CapaciteSavoirfaire.capaciteSavoirFaire" ----
# (moved down from row 3, content unchanged)
$ws.Range("A6").Value = "CapaciteSavoirfaire.capaciteSavoirFaire"
$ws.Range("B6").Value = "CapaciteSavoirfaire.capaciteSavoirFaire"
$ws.Range("F6").Value = "0"
$ws.Range("G6").Value = "1"
$ws.Range("K6").Value = "Coding`n"
$ws.Range("L6").Value = " Savoir-faire de type capacit" + [char]0x00E9 + " de m" + [char]0x00E9 + "decine (dipl" + [char]0x00F4 + "me permettant aux docteurs en m" + [char]0x00E9 + "decine d'acqu" + [char]0x00E9 + "rir des comp" + [char]0x00E9 + "tences suppl" + [char]0x00E9 + "mentaires)."
$ws.Range("M6").Value = $ws.Range("L6").Value
$ws.Range("X6").Value = "preferred"
$ws.Range("Y6").Value = ""
$ws.Range("Z6").Value = "https://interop.esante.gouv.fr/ig/fhir/mos/ValueSet/capaciteSavoirFaire-vs"
$ws.Range("AF6").Value = "CapaciteSavoirfaire.capaciteSavoirFaire"
$ws.Range("AG6").Value = "0"
$ws.Range("AH6").Value = "1"

# ---------------------------------------------------------------------
# 3) Column width tweaks (A, B widen slightly; Z widens for the longer
#    binding value-set URL now shown).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 32.1666666666667
$ws.Columns.Item(2).ColumnWidth = 32.1666666666667
$ws.Columns.Item(26).ColumnWidth = 68.0
